# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "September Sun" durazno at the top
# of the Vega Monumental Concepción block (old row 325), pushing the
# existing rows 325-358 down to 327-360.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("325:326").Insert()

# Row 325: September Sun / Especial
$ws.Range("A325").Value = 11
$ws.Range("B325").Value = "Vega Monumental Concepción"
$ws.Range("C325").Value = "Bíobío"
$ws.Range("D325").Value = 44995
$ws.Range("E325").Value = 8
$ws.Range("F325").Value = "Fruta"
$ws.Range("G325").Value = 100103
$ws.Range("H325").Value = "Frutos de hueso (carozo)"
$ws.Range("I325").Value = 100103004
$ws.Range("J325").Value = "Durazno"
$ws.Range("K325").Value = "September Sun"
$ws.Range("L325").Value = "Especial"
$ws.Range("M325").Value = 250
$ws.Range("N325").Value = 16000
$ws.Range("O325").Value = 16000
$ws.Range("P325").Value = 16000
$ws.Range("Q325").Value = "$/caja 16 kilos empedrada"
$ws.Range("R325").Value = "Región de O'Higgins"
$ws.Range("S325").Value = 1000
$ws.Range("T325").Value = 16

# Row 326: September Sun / Primera
$ws.Range("A326").Value = 11
$ws.Range("B326").Value = "Vega Monumental Concepción"
$ws.Range("C326").Value = "Bíobío"
$ws.Range("D326").Value = 44995
$ws.Range("E326").Value = 8
$ws.Range("F326").Value = "Fruta"
$ws.Range("G326").Value = 100103
$ws.Range("H326").Value = "Frutos de hueso (carozo)"
$ws.Range("I326").Value = 100103004
$ws.Range("J326").Value = "Durazno"
$ws.Range("K326").Value = "September Sun"
$ws.Range("L326").Value = "Primera"
$ws.Range("M326").Value = 300
$ws.Range("N326").Value = 14000
$ws.Range("O326").Value = 14000
$ws.Range("P326").Value = 14000
$ws.Range("Q326").Value = "$/caja 16 kilos empedrada"
$ws.Range("R326").Value = "Región de O'Higgins"
$ws.Range("S326").Value = 875
$ws.Range("T326").Value = 16
